# Update "a_prime_leet" Problem sheet with round 2 (PSAT prime mock exam) data.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Problem")

# --- Correct the existing round-1 answer values (column E, rows 2-71) ---
$ws.Range("E2").Value = 2
$ws.Range("E3").Value = 5
$ws.Range("E4").Value = 5
$ws.Range("E5").Value = 2
$ws.Range("E6").Value = 3
$ws.Range("E8").Value = 2
$ws.Range("E9").Value = 1
$ws.Range("E11").Value = 3
$ws.Range("E12").Value = 4
$ws.Range("E13").Value = 2
$ws.Range("E14").Value = 3
$ws.Range("E15").Value = 5
$ws.Range("E16").Value = 5
$ws.Range("E17").Value = 5
$ws.Range("E18").Value = 2
$ws.Range("E19").Value = 3
$ws.Range("E20").Value = 3
$ws.Range("E22").Value = 3
$ws.Range("E23").Value = 2
$ws.Range("E24").Value = 4
$ws.Range("E26").Value = 5
$ws.Range("E28").Value = 2
$ws.Range("E30").Value = 2
$ws.Range("E31").Value = 1
$ws.Range("E32").Value = 3
$ws.Range("E34").Value = 3
$ws.Range("E36").Value = 3
$ws.Range("E37").Value = 1
$ws.Range("E39").Value = 1
$ws.Range("E40").Value = 1
$ws.Range("E41").Value = 2
$ws.Range("E42").Value = 3
$ws.Range("E43").Value = 3
$ws.Range("E44").Value = 1
$ws.Range("E45").Value = 4
$ws.Range("E46").Value = 2
$ws.Range("E47").Value = 3
$ws.Range("E48").Value = 4
$ws.Range("E49").Value = 4
$ws.Range("E50").Value = 5
$ws.Range("E53").Value = 5
$ws.Range("E54").Value = 3
$ws.Range("E55").Value = 2
$ws.Range("E56").Value = 3
$ws.Range("E57").Value = 2
$ws.Range("E58").Value = 3
$ws.Range("E60").Value = 1
$ws.Range("E61").Value = 1
$ws.Range("E62").Value = 4
$ws.Range("E63").Value = 5
$ws.Range("E64").Value = 5
$ws.Range("E65").Value = 5
$ws.Range("E66").Value = 5
$ws.Range("E68").Value = 2
$ws.Range("E70").Value = 2
$ws.Range("E71").Value = 1

# --- Append round-2 data (rows 72-106): id, leet_id, subject, number, answer ---
$ws.Range("B72").Value = 2
$ws.Range("C72").Value = "언어"
$ws.Range("D72").Value = 1
$ws.Range("E72").Value = 2
$ws.Range("B73").Value = 2
$ws.Range("C73").Value = "언어"
$ws.Range("D73").Value = 2
$ws.Range("E73").Value = 5
$ws.Range("B74").Value = 2
$ws.Range("C74").Value = "언어"
$ws.Range("D74").Value = 3
$ws.Range("E74").Value = 5
$ws.Range("B75").Value = 2
$ws.Range("C75").Value = "언어"
$ws.Range("D75").Value = 4
$ws.Range("E75").Value = 2
$ws.Range("B76").Value = 2
$ws.Range("C76").Value = "언어"
$ws.Range("D76").Value = 5
$ws.Range("E76").Value = 3
$ws.Range("B77").Value = 2
$ws.Range("C77").Value = "언어"
$ws.Range("D77").Value = 6
$ws.Range("E77").Value = 5
$ws.Range("B78").Value = 2
$ws.Range("C78").Value = "언어"
$ws.Range("D78").Value = 7
$ws.Range("E78").Value = 2
$ws.Range("B79").Value = 2
$ws.Range("C79").Value = "언어"
$ws.Range("D79").Value = 8
$ws.Range("E79").Value = 1
$ws.Range("B80").Value = 2
$ws.Range("C80").Value = "언어"
$ws.Range("D80").Value = 9
$ws.Range("E80").Value = 2
$ws.Range("B81").Value = 2
$ws.Range("C81").Value = "언어"
$ws.Range("D81").Value = 10
$ws.Range("E81").Value = 3
$ws.Range("B82").Value = 2
$ws.Range("C82").Value = "언어"
$ws.Range("D82").Value = 11
$ws.Range("E82").Value = 1
$ws.Range("B83").Value = 2
$ws.Range("C83").Value = "언어"
$ws.Range("D83").Value = 12
$ws.Range("E83").Value = 1
$ws.Range("B84").Value = 2
$ws.Range("C84").Value = "언어"
$ws.Range("D84").Value = 13
$ws.Range("E84").Value = 5
$ws.Range("B85").Value = 2
$ws.Range("C85").Value = "언어"
$ws.Range("D85").Value = 14
$ws.Range("E85").Value = 4
$ws.Range("B86").Value = 2
$ws.Range("C86").Value = "언어"
$ws.Range("D86").Value = 15
$ws.Range("E86").Value = 4
$ws.Range("B87").Value = 2
$ws.Range("C87").Value = "추리"
$ws.Range("D87").Value = 1
$ws.Range("E87").Value = 3
$ws.Range("B88").Value = 2
$ws.Range("C88").Value = "추리"
$ws.Range("D88").Value = 2
$ws.Range("E88").Value = 3
$ws.Range("B89").Value = 2
$ws.Range("C89").Value = "추리"
$ws.Range("D89").Value = 3
$ws.Range("E89").Value = 1
$ws.Range("B90").Value = 2
$ws.Range("C90").Value = "추리"
$ws.Range("D90").Value = 4
$ws.Range("E90").Value = 3
$ws.Range("B91").Value = 2
$ws.Range("C91").Value = "추리"
$ws.Range("D91").Value = 5
$ws.Range("E91").Value = 1
$ws.Range("B92").Value = 2
$ws.Range("C92").Value = "추리"
$ws.Range("D92").Value = 6
$ws.Range("E92").Value = 1
$ws.Range("B93").Value = 2
$ws.Range("C93").Value = "추리"
$ws.Range("D93").Value = 7
$ws.Range("E93").Value = 3
$ws.Range("B94").Value = 2
$ws.Range("C94").Value = "추리"
$ws.Range("D94").Value = 8
$ws.Range("E94").Value = 1
$ws.Range("B95").Value = 2
$ws.Range("C95").Value = "추리"
$ws.Range("D95").Value = 9
$ws.Range("E95").Value = 4
$ws.Range("B96").Value = 2
$ws.Range("C96").Value = "추리"
$ws.Range("D96").Value = 10
$ws.Range("E96").Value = 2
$ws.Range("B97").Value = 2
$ws.Range("C97").Value = "추리"
$ws.Range("D97").Value = 11
$ws.Range("E97").Value = 4
$ws.Range("B98").Value = 2
$ws.Range("C98").Value = "추리"
$ws.Range("D98").Value = 12
$ws.Range("E98").Value = 5
$ws.Range("B99").Value = 2
$ws.Range("C99").Value = "추리"
$ws.Range("D99").Value = 13
$ws.Range("E99").Value = 3
$ws.Range("B100").Value = 2
$ws.Range("C100").Value = "추리"
$ws.Range("D100").Value = 14
$ws.Range("E100").Value = 3
$ws.Range("B101").Value = 2
$ws.Range("C101").Value = "추리"
$ws.Range("D101").Value = 15
$ws.Range("E101").Value = 2
$ws.Range("B102").Value = 2
$ws.Range("C102").Value = "추리"
$ws.Range("D102").Value = 16
$ws.Range("E102").Value = 1
$ws.Range("B103").Value = 2
$ws.Range("C103").Value = "추리"
$ws.Range("D103").Value = 17
$ws.Range("E103").Value = 4
$ws.Range("B104").Value = 2
$ws.Range("C104").Value = "추리"
$ws.Range("D104").Value = 18
$ws.Range("E104").Value = 5
$ws.Range("B105").Value = 2
$ws.Range("C105").Value = "추리"
$ws.Range("D105").Value = 19
$ws.Range("E105").Value = 2
$ws.Range("B106").Value = 2
$ws.Range("C106").Value = "추리"
$ws.Range("D106").Value = 20
$ws.Range("E106").Value = 1

# --- Re-enter the CSE array formula over the new full range so it recalculates/fills A2:A106 ---
$ws.Range("A2:A106").FormulaArray = "=SEQUENCE(COUNT(B:B))"

# --- Update the view state to match where the user ended up after entering the new rows ---
[void]$ws.Range("E106").Select()
